$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.726.65'
$ws.Range("E2").Value = '  -4.09%  '
$ws.Range("D3").Value = '1.818.02'
$ws.Range("E3").Value = '  -2.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.95'
$ws.Range("E5").Value = '  -7.79%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  -5.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3533'
$ws.Range("E8").Value = '  -5.66%  '
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06669'
$ws.Range("E10").Value = '  -7.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.08'
$ws.Range("E11").Value = '  -7.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8282'
$ws.Range("E12").Value = '  -6.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07910'
$ws.Range("D14").Value = '1.804.44'
$ws.Range("E14").Value = '  -3.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.077'
$ws.Range("E15").Value = '  -4.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.71'
$ws.Range("E16").Value = '  -6.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.09'
$ws.Range("E18").Value = '  -5.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008028'
$ws.Range("E19").Value = '  -5.94%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '25.767.56'
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.744'
$ws.Range("E22").Value = '  -4.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.996'
$ws.Range("E23").Value = '  -5.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.107'
$ws.Range("E24").Value = '  -4.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.224'
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.81'
$ws.Range("E26").Value = '  -3.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.672'
$ws.Range("E27").Value = '  -3.48%  '
$ws.Range("E28").Value = '  -5.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.33'
$ws.Range("E29").Value = '  -4.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.321'
$ws.Range("E30").Value = '  -8.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.240'
$ws.Range("E31").Value = '  -8.20%  '
$ws.Range("E32").Value = '  -4.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04887'
$ws.Range("E33").Value = '  -2.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7286'
$ws.Range("E34").Value = '  -10.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.861'
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.379'
$ws.Range("E39").Value = '  -9.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01855'
$ws.Range("E40").Value = '  -5.27%  '
$ws.Range("E41").Value = '  -14.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9667'
$ws.Range("E42").Value = '  -9.58%  '
$ws.Range("E43").Value = '  -6.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.21'
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.031'
$ws.Range("E45").Value = '  -9.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4565'
$ws.Range("E47").Value = '  -10.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1365'
$ws.Range("E48").Value = '  -8.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.51'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.206'
$ws.Range("E50").Value = '  -7.50%  '
$ws.Range("E51").Value = '  -8.74%  '
